# Results_1-3.xlsx -- "various updates on Part 1.3 and Part 1.4"
#
# Summary of the edit:
#  - I1 header text "Ration" was a typo; fixed to "Ratio".
#  - The "Ratio" (I column, Restaurant) computation is now mirrored into a new
#    L column ("Ratio" for Airline), with a header label in L1.
#  - A1 gets a blank placeholder label (single space) to match the new row.
#  - The Ratio result columns (I2:I5 and L2:L5) get a "0.0" number format.
#  - Column L is sized to fit its new header/content.
#  - The active selection moves to A9 (cursor location after the edit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text / shared-string content -----------------------------------------

# Fix the "Ration" -> "Ratio" typo in the existing header (I1). Doing this
# first means the corrected "Ratio" string is appended to the shared string
# table right after the now-orphaned "Ration" entry is dropped.
$ws.Range("I1").Value = "Ratio"

# New placeholder cell at A1 (single space), matching the new row 1 layout.
$ws.Range("A1").Value = " "

# New column L mirrors column I: same header text and formatting.
$ws.Range("L1").Value = "Ratio"
$ws.Range("L1").WrapText = $true

# Re-index the "Bag of Words Count for ..." labels stay the same text (no
# content change needed -- only their shared-string index shifts, which
# happens automatically once "Ration" is removed above).
$ws.Range("A8").Value = "Bag of Words Count for Airline "
$ws.Range("A11").Value = "Bag of Words Count for Restaurant"

# --- Number formatting -----------------------------------------------------

$ws.Range("I2:I5").NumberFormat = "0.0"
$ws.Range("L2:L5").NumberFormat = "0.0"

# --- Column sizing -----------------------------------------------------

$ws.Columns.Item(12).ColumnWidth = 10.8

# --- Selection ---------------------------------------------------------

$ws.Range("A9").Select()
